# Created experiment order generation script
# Re-generates the per-task-order sheets for participant_49 with a fresh
# randomization pass: sheet tab names and their "go/GNG / OB/TB/ZB / eyes /
# MM/ZM / vSAT/SAT" stim-file contents are regenerated for a new run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (was GNG_TO-...) -> vSAT task order
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515890482248614"
$ws1.Cells.Item(2,2).Value = "vSAT_stims-16515890482092366.csv"
$ws1.Cells.Item(3,2).Value = "SAT_stims-1651589048162331.csv"
$ws1.Cells.Item(4,2).Value = "SAT_stims-16515890481779914.csv"
$ws1.Cells.Item(5,2).Value = "vSAT_stims-16515890481936162.csv"

# ---------------------------------------------------------------------
# Sheet 2 (was NB_TO-...) -> NB task order (stays NB, new file names)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16515890493298016"
$ws2.Cells.Item(2,2).Value = "OB-16515890489335659.csv"
$ws2.Cells.Item(3,2).Value = "TB-16515890493141346.csv"
$ws2.Cells.Item(4,2).Value = "ZB-match_7-16515890484605336.csv"
$ws2.Cells.Item(5,2).Value = "TB-16515890492169983.csv"
$ws2.Cells.Item(6,2).Value = "ZB-match_3-16515890486501694.csv"
$ws2.Cells.Item(7,2).Value = "ZB-match_5-1651589048346677.csv"
$ws2.Cells.Item(8,2).Value = "OB-16515890490273166.csv"
$ws2.Cells.Item(9,2).Value = "OB-16515890486970026.csv"
$ws2.Cells.Item(10,2).Value = "TB-16515890491232517.csv"

# ---------------------------------------------------------------------
# Sheet 3 (was RS_TO-..., eyes open/closed, A1:B3) -> TOL task order
# (MM_stims / ZM_stims, A1:B7) -- grows from 2 data rows to 6 data rows.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "TOL_TO-16515890493789148"

# Extend the styled index column (A) down to row 7, copying the existing
# formatted cell so the new rows keep the same bold/border/centered style.
$ws3.Range("A3").Copy($ws3.Range("A4:A7"))

$ws3.Cells.Item(2,1).Value = 0
$ws3.Cells.Item(2,2).Value = "MM_stims-16515890493453526.csv"
$ws3.Cells.Item(3,1).Value = 1
$ws3.Cells.Item(3,2).Value = "ZM_stims-16515890493298016.csv"
$ws3.Cells.Item(4,1).Value = 2
$ws3.Cells.Item(4,2).Value = "MM_stims-16515890493632903.csv"
$ws3.Cells.Item(5,1).Value = 3
$ws3.Cells.Item(5,2).Value = "ZM_stims-16515890493464952.csv"
$ws3.Cells.Item(6,1).Value = 4
$ws3.Cells.Item(6,2).Value = "MM_stims-16515890493789148.csv"
$ws3.Cells.Item(7,1).Value = 5
$ws3.Cells.Item(7,2).Value = "ZM_stims-16515890493632903.csv"

# ---------------------------------------------------------------------
# Sheet 4 (was TOL_TO-..., MM_stims/ZM_stims, A1:B7) -> RS task order
# (eyes open/closed, A1:B3) -- shrinks from 6 data rows to 2 data rows.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "RS_TO-16515890493789148"

# Drop the now-unused trailing rows (4 through 7) so the sheet shrinks
# back down to the A1:B3 extent used by the RS task order.
$ws4.Range("A4:B7").EntireRow.Delete()

$ws4.Cells.Item(2,2).Value = "eyes closed"
$ws4.Cells.Item(3,2).Value = "eyes open"

# ---------------------------------------------------------------------
# Sheet 5 (was vSAT_TO-...) -> GNG task order
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "GNG_TO-16515890494101632"
$ws5.Cells.Item(2,2).Value = "go_stims-16515890493789148.csv"
$ws5.Cells.Item(3,2).Value = "GNG_stims-16515890493945398.csv"
$ws5.Cells.Item(4,2).Value = "go_stims-16515890493945398.csv"
$ws5.Cells.Item(5,2).Value = "GNG_stims-16515890494101632.csv"
